$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9019181025842329
$ws.Range("D2").Value = 0.2950447304472732
$ws.Range("E2").Value = 0.3159889265604434
$ws.Range("F2").Value = 1.49810685513463
$ws.Range("G2").Value = 0.780123361248485
$ws.Range("H2").Value = 0.8762598520114366
$ws.Range("J2").Value = 0.4149153750845276
$ws.Range("K2").Value = 0.3636192724592604
$ws.Range("L2").Value = 0.09193412346857599
$ws.Range("M2").Value = 0.1723180485962175
$ws.Range("O2").Value = 3.326857037664681
$ws.Range("B3").Value = 0.8797045962019183
$ws.Range("D3").Value = 0.2955555724836358
$ws.Range("E3").Value = 0.3182675989860906
$ws.Range("F3").Value = 1.507518249656592
$ws.Range("G3").Value = 0.7854615810608649
$ws.Range("H3").Value = 0.8824372293548208
$ws.Range("J3").Value = 0.4178062770101558
$ws.Range("K3").Value = 0.3202758800964034
$ws.Range("L3").Value = 0.08478178247589341
$ws.Range("M3").Value = 0.1656910795737474
$ws.Range("O3").Value = 3.350857000625552
$ws.Range("B4").Value = 0.866378320305472
$ws.Range("D4").Value = 0.2959638960239417
$ws.Range("E4").Value = 0.3197544185701204
$ws.Range("F4").Value = 1.513971668206437
$ws.Range("G4").Value = 0.7891719310338559
$ws.Range("H4").Value = 0.8865560561696668
$ws.Range("J4").Value = 0.4196805312810259
$ws.Range("K4").Value = 0.2935319630387028
$ws.Range("L4").Value = 0.08039613549231461
$ws.Range("M4").Value = 0.1616700459861136
$ws.Range("O4").Value = 3.367183790651126
$ws.Range("B5").Value = 0.8610270484147122
$ws.Range("D5").Value = 0.2961541678882682
$ws.Range("E5").Value = 0.3203824022191943
$ws.Range("F5").Value = 1.516771350658146
$ws.Range("G5").Value = 0.7907927616175883
$ws.Range("H5").Value = 0.8883165577976584
$ws.Range("J5").Value = 0.4204693016269978
$ws.Range("K5").Value = 0.2826014131433965
$ws.Range("L5").Value = 0.07861058082998795
$ws.Range("M5").Value = 0.1600436546131903
$ws.Range("O5").Value = 3.374237365397263
$ws.Range("B6").Value = 0.8601432809236371
$ws.Range("D6").Value = 0.2961872066136948
$ws.Range("E6").Value = 0.3204880138910475
$ws.Range("F6").Value = 1.517246500386598
$ws.Range("G6").Value = 0.7910684736088029
$ws.Range("H6").Value = 0.8886138460258408
$ws.Range("J6").Value = 0.4206017873839381
$ws.Range("K6").Value = 0.2807844811798645
$ws.Range("L6").Value = 0.0783141936460936
$ws.Range("M6").Value = 0.1597743357182644
$ws.Range("O6").Value = 3.375432789880577
$ws.Range("B7").Value = 0.8663058293124948
$ws.Range("D7").Value = 0.2959663653242757
$ws.Range("E7").Value = 0.3197627982738984
$ws.Range("F7").Value = 1.514008737696571
$ws.Range("G7").Value = 0.7891933494135301
$ws.Range("H7").Value = 0.8865794665779703
$ws.Range("J7").Value = 0.4196910676443384
$ws.Range("K7").Value = 0.2933846791738688
$ws.Range("L7").Value = 0.08037204807855147
$ws.Range("M7").Value = 0.1616480622502614
$ws.Range("O7").Value = 3.367277296560346
$ws.Range("B8").Value = 0.8941942993738792
$ws.Range("D8").Value = 0.2952012643208235
$ws.Range("E8").Value = 0.3167564327122987
$ws.Range("F8").Value = 1.501211961057926
$ws.Range("G8").Value = 0.7818742020421396
$ws.Range("H8").Value = 0.8783222316777781
$ws.Range("J8").Value = 0.4158915901488367
$ws.Range("K8").Value = 0.3487020626762956
$ws.Range("L8").Value = 0.08946685155007117
$ws.Range("M8").Value = 0.1700232074918837
$ws.Range("O8").Value = 3.334802237024178
$ws.Range("B9").Value = 0.951342859104841
$ws.Range("D9").Value = 0.294448987035075
$ws.Range("E9").Value = 0.3115551057863071
$ws.Range("F9").Value = 1.481464051879975
$ws.Range("G9").Value = 0.7709531031455157
$ws.Range("H9").Value = 0.8647113359931495
$ws.Range("J9").Value = 0.4092260884239867
$ws.Range("K9").Value = 0.4561136586841883
$ws.Range("L9").Value = 0.107343277193678
$ws.Range("M9").Value = 0.1868214497294893
$ws.Range("O9").Value = 3.283729004185773
$ws.Range("B10").Value = 0.9948027389778247
$ws.Range("D10").Value = 0.2943484805699939
$ws.Range("E10").Value = 0.3081543913434563
$ws.Range("F10").Value = 1.470205079164494
$ws.Range("G10").Value = 0.7650203041600321
$ws.Range("H10").Value = 0.8562794270650329
$ws.Range("J10").Value = 0.404804704786919
$ws.Range("K10").Value = 0.5343504363836757
$ws.Range("L10").Value = 0.1204964295988304
$ws.Range("M10").Value = 0.199384946888884
$ws.Range("O10").Value = 3.253878987043123
$ws.Range("B11").Value = 1.014888281255082
$ws.Range("D11").Value = 0.294400143990643
$ws.Range("E11").Value = 0.3066981432947387
$ws.Range("F11").Value = 1.465786767690886
$ws.Range("G11").Value = 0.762775141620466
$ws.Range("H11").Value = 0.8527828088493123
$ws.Range("J11").Value = 0.4028959801192258
$ws.Range("K11").Value = 0.5697893885626399
$ws.Range("L11").Value = 0.1264831533836173
$ws.Range("M11").Value = 0.205147272631848
$ws.Range("O11").Value = 3.241962791590538
$ws.Range("B12").Value = 1.022538959998201
$ws.Range("D12").Value = 0.294433641149773
$ws.Range("E12").Value = 0.3061597111347885
$ws.Range("F12").Value = 1.46421465496725
$ws.Range("G12").Value = 0.7619901741382904
$ws.Range("H12").Value = 0.8515073983572279
$ws.Range("J12").Value = 0.4021879016942691
$ws.Range("K12").Value = 0.5831867950992091
$ws.Range("L12").Value = 0.128750514637801
$ws.Range("M12").Value = 0.2073359538650621
$ws.Range("O12").Value = 3.237689294514354
$ws.Range("B13").Value = 1.020889270666117
$ws.Range("D13").Value = 0.2944258083682598
$ws.Range("E13").Value = 0.3062750937059988
$ws.Range("F13").Value = 1.46454874769136
$ws.Range("G13").Value = 0.7621563301852348
$ws.Range("H13").Value = 0.8517799168976268
$ws.Range("J13").Value = 0.4023397452268362
$ws.Range("K13").Value = 0.5803024396256546
$ws.Range("L13").Value = 0.1282621860889606
$ws.Range("M13").Value = 0.2068642906607323
$ws.Range("O13").Value = 3.238599046144117
$ws.Range("B14").Value = 1.015516814599948
$ws.Range("D14").Value = 0.2944026209798878
$ws.Range("E14").Value = 0.3066535854574468
$ws.Range("F14").Value = 1.465655405579682
$ws.Range("G14").Value = 0.7627092547148209
$ws.Range("H14").Value = 0.8526769048345528
$ws.Range("J14").Value = 0.4028374313999583
$ws.Range("K14").Value = 0.5708920578314576
$ws.Range("L14").Value = 0.1266696850483555
$ws.Range("M14").Value = 0.2053272051438668
$ws.Range("O14").Value = 3.241606421447244
$ws.Range("B15").Value = 1.012231835052916
$ws.Range("D15").Value = 0.2943902304793156
$ws.Range("E15").Value = 0.3068871168552825
$ws.Range("F15").Value = 1.46634641499935
$ws.Range("G15").Value = 0.7630564309774286
$ws.Range("H15").Value = 0.8532326734571285
$ws.Range("J15").Value = 0.4031441936683819
$ws.Range("K15").Value = 0.5651249664263105
$ws.Range("L15").Value = 0.1256942691425991
$ws.Range("M15").Value = 0.2043865523156398
$ws.Range("O15").Value = 3.243479632792287
$ws.Range("B16").Value = 0.9934963941563524
$ws.Range("D16").Value = 0.2943470578438934
$ws.Range("E16").Value = 0.3082513839607133
$ws.Range("F16").Value = 1.470507968682583
$ws.Range("G16").Value = 0.7651761591605606
$ws.Range("H16").Value = 0.8565147560647546
$ws.Range("J16").Value = 0.404931505497709
$ws.Range("K16").Value = 0.5320313042271607
$ws.Range("L16").Value = 0.1201052361214039
$ws.Range("M16").Value = 0.1990092997081661
$ws.Range("O16").Value = 3.254691180492784
$ws.Range("B17").Value = 0.9820831417023328
$ws.Range("D17").Value = 0.2943454651641417
$ws.Range("E17").Value = 0.3091115384216545
$ws.Range("F17").Value = 1.473241011656476
$ws.Range("G17").Value = 0.7665927366237639
$ws.Range("H17").Value = 0.8586150025407733
$ws.Range("J17").Value = 0.4060542111400656
$ws.Range("K17").Value = 0.5116900857727558
$ws.Range("L17").Value = 0.1166772779876055
$ws.Range("M17").Value = 0.1957224836184395
$ws.Range("O17").Value = 3.261994830899368
$ws.Range("B18").Value = 0.9755482681241858
$ws.Range("D18").Value = 0.2943537195734578
$ws.Range("E18").Value = 0.3096148206733362
$ws.Range("F18").Value = 1.474879206860209
$ws.Range("G18").Value = 0.7674502195709962
$ws.Range("H18").Value = 0.8598549311235644
$ws.Range("J18").Value = 0.4067096197222666
$ws.Range("K18").Value = 0.4999761562296214
$ws.Range("L18").Value = 0.1147059258708651
$ws.Range("M18").Value = 0.1938364395822418
$ws.Range("O18").Value = 3.266352212643554
$ws.Range("B19").Value = 0.973340797371435
$ws.Range("D19").Value = 0.2943580915590047
$ws.Range("E19").Value = 0.3097866918698902
$ws.Range("F19").Value = 1.475445249708507
$ws.Range("G19").Value = 0.7677478836220075
$ws.Range("H19").Value = 0.8602802347493537
$ws.Range("J19").Value = 0.4069331897224782
$ws.Range("K19").Value = 0.4960076064868417
$ws.Range("L19").Value = 0.1140385193607756
$ws.Range("M19").Value = 0.193198626006172
$ws.Range("O19").Value = 3.267854434661302
$ws.Range("B20").Value = 0.9832950298639389
$ws.Range("D20").Value = 0.2943446861085661
$ws.Range("E20").Value = 0.3090190894623639
$ws.Range("F20").Value = 1.472943222005526
$ws.Range("G20").Value = 0.7664375198224604
$ws.Range("H20").Value = 0.8583881244475293
$ws.Range("J20").Value = 0.405933697914227
$ws.Range("K20").Value = 0.5138569184570656
$ws.Range("L20").Value = 0.1170421581121275
$ws.Range("M20").Value = 0.1960719118793577
$ws.Range("O20").Value = 3.261201147485224
$ws.Range("B21").Value = 1.017093626573228
$ws.Range("D21").Value = 0.2944090540913891
$ws.Range("E21").Value = 0.3065420602292539
$ws.Range("F21").Value = 1.465327613254104
$ws.Range("G21").Value = 0.7625450772574709
$ws.Range("H21").Value = 0.8524121171229382
$ws.Range("J21").Value = 0.4026908498645883
$ws.Range("K21").Value = 0.5736567335193854
$ws.Range("L21").Value = 0.1271374337779463
$ws.Range("M21").Value = 0.2057785062309634
$ws.Range("O21").Value = 3.240716600063962
$ws.Range("B22").Value = 1.039443329149663
$ws.Range("D22").Value = 0.2945323038761387
$ws.Range("E22").Value = 0.3049990384745698
$ws.Range("F22").Value = 1.46093904289453
$ws.Range("G22").Value = 0.7603813258043601
$ws.Range("H22").Value = 0.8487901757608967
$ws.Range("J22").Value = 0.4006572096198528
$ws.Range("K22").Value = 0.6126075334467203
$ws.Range("L22").Value = 0.1337370675456953
$ws.Range("M22").Value = 0.2121607998864974
$ws.Range("O22").Value = 3.228721210778872
$ws.Range("B23").Value = 1.027491259314445
$ws.Range("D23").Value = 0.2944591182606189
$ws.Range("E23").Value = 0.3058156478638709
$ws.Range("F23").Value = 1.4632274916558
$ws.Range("G23").Value = 0.7615013786724631
$ws.Range("H23").Value = 0.8506973392669863
$ws.Range("J23").Value = 0.401734767967346
$ws.Range("K23").Value = 0.5918310822531794
$ws.Range("L23").Value = 0.1302146060222498
$ws.Range("M23").Value = 0.2087509832602734
$ws.Range("O23").Value = 3.234996030110295
$ws.Range("B24").Value = 0.9827470518326891
$ws.Range("D24").Value = 0.2943450097515736
$ws.Range("E24").Value = 0.3090608583232966
$ws.Range("F24").Value = 1.47307764420713
$ws.Range("G24").Value = 0.7665075591662571
$ws.Range("H24").Value = 0.8584905948891333
$ws.Range("J24").Value = 0.405988150946726
$ws.Range("K24").Value = 0.5128773540569966
$ws.Range("L24").Value = 0.1168771975824257
$ws.Range("M24").Value = 0.195913924167094
$ws.Range("O24").Value = 3.261559478289229
$ws.Range("B25").Value = 0.9356221077760267
$ws.Range("D25").Value = 0.2945728099504166
$ws.Range("E25").Value = 0.3128881453681593
$ws.Range("F25").Value = 1.486234990140375
$ws.Range("G25").Value = 0.7735402753325289
$ws.Range("H25").Value = 0.8681176699601778
$ws.Range("J25").Value = 0.4109455225704739
$ws.Range("K25").Value = 0.4271730568327143
$ws.Range("L25").Value = 0.1025034032092265
$ws.Range("M25").Value = 0.1822376534690022
$ws.Range("O25").Value = 3.296197103262173
